$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numRows = 8
$numCols = 8
$startRow = 534

$arr = New-Object 'object[,]' $numRows,$numCols

$arr[0,0] = 'linearization_heuristic_optBouncing=False_initial_uhat=age_group_gradient_targetGroups=True_targetAct=False_targetTests=False'
$arr[0,1] = 50
$arr[0,2] = 2900
$arr[0,3] = 60000
$arr[0,4] = 0.1
$arr[0,5] = 52.97564769386039
$arr[0,6] = 1868.070516240173
$arr[0,7] = 49501127135.07371
$arr[1,0] = 'linearization_heuristic_optBouncing=False_initial_uhat=age_group_gradient_targetGroups=True_targetAct=False_targetTests=False'
$arr[1,1] = 50
$arr[1,2] = 2900
$arr[1,3] = 120000
$arr[1,4] = 0.1
$arr[1,5] = 53.06393615750851
$arr[1,6] = 1821.193960486289
$arr[1,7] = 49676603718.91109
$arr[2,0] = 'linearization_heuristic_optBouncing=False_initial_uhat=age_group_gradient_targetGroups=True_targetAct=False_targetTests=False'
$arr[2,1] = 50
$arr[2,2] = 2900
$arr[2,3] = 240000
$arr[2,4] = 0.1
$arr[2,5] = 52.69334229288796
$arr[2,6] = 1426.416571816286
$arr[2,7] = 50040276650.5134
$arr[3,0] = 'linearization_heuristic_optBouncing=False_initial_uhat=age_group_gradient_targetGroups=True_targetAct=False_targetTests=False'
$arr[3,1] = 50
$arr[3,2] = 2900
$arr[3,3] = 360000
$arr[3,4] = 0.1
$arr[3,5] = 52.70093459161686
$arr[3,6] = 1801.951147482135
$arr[3,7] = 49349392851.93075
$arr[4,0] = 'linearization_heuristic_optBouncing=False_initial_uhat=age_group_gradient_targetGroups=True_targetAct=False_targetTests=True'
$arr[4,1] = 50
$arr[4,2] = 2900
$arr[4,3] = 60000
$arr[4,4] = 0.1
$arr[4,5] = 52.98238601317134
$arr[4,6] = 1776.946010866787
$arr[4,7] = 49677352614.84064
$arr[5,0] = 'linearization_heuristic_optBouncing=False_initial_uhat=age_group_gradient_targetGroups=True_targetAct=False_targetTests=True'
$arr[5,1] = 50
$arr[5,2] = 2900
$arr[5,3] = 120000
$arr[5,4] = 0.1
$arr[5,5] = 53.07388983985473
$arr[5,6] = 1680.452738282573
$arr[5,7] = 49948329248.60695
$arr[6,0] = 'linearization_heuristic_optBouncing=False_initial_uhat=age_group_gradient_targetGroups=True_targetAct=False_targetTests=True'
$arr[6,1] = 50
$arr[6,2] = 2900
$arr[6,3] = 240000
$arr[6,4] = 0.1
$arr[6,5] = 52.97304648324432
$arr[6,6] = 1411.228643471699
$arr[6,7] = 50348229650.97617
$arr[7,0] = 'linearization_heuristic_optBouncing=False_initial_uhat=age_group_gradient_targetGroups=True_targetAct=False_targetTests=True'
$arr[7,1] = 50
$arr[7,2] = 2900
$arr[7,3] = 360000
$arr[7,4] = 0.1
$arr[7,5] = 52.71130665491413
$arr[7,6] = 1579.370988188422
$arr[7,7] = 49773753216.37659

$targetRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($startRow + $numRows - 1, $numCols))
$targetRange.Value2 = $arr

Write-Host "Wrote $numRows rows starting at row $startRow"
